$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 165, shifting existing rows 165-176 down to 166-177.
$ws.Rows.Item(165).Insert()

# Fill the new row 165 with a copy of the (now shifted) row 166 values, then
# update the cells that actually differ for the new weekly entry.
$ws.Range("A166:R166").Copy($ws.Range("A165:R165"))

$ws.Range("D165").Value = 44610
$ws.Range("J165").Value = 60
$ws.Range("K165").Value = 8000
$ws.Range("L165").Value = 8500
$ws.Range("M165").Value = 8250
$ws.Range("O165").Value = "Provincia del Elquí"
$ws.Range("P165").Value = 1375
